$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert a new data row right after the existing row 112
# (row 113), pushing the previously-existing rows 113:146 down to
# 114:147 unchanged.
$ws.Rows.Item(113).Insert()

$ws.Cells.Item(113, 1).Value = 8
$ws.Cells.Item(113, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(113, 3).Value = "Coquimbo"
$ws.Cells.Item(113, 4).Value = 44559
$ws.Cells.Item(113, 5).Value = 4
$ws.Cells.Item(113, 6).Value = 100112037
$ws.Cells.Item(113, 7).Value = "Cebollín"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 3000
$ws.Cells.Item(113, 11).Value = 900
$ws.Cells.Item(113, 12).Value = 1000
$ws.Cells.Item(113, 13).Value = 950
$ws.Cells.Item(113, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(113, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(113, 16).Value = 158
$ws.Cells.Item(113, 17).Value = 6
$ws.Cells.Item(113, 18).Value = "Hortaliza"
